# Remove the trailing "empty line / Ver no Jupiter.../ (c) 2020 ..." footer
# paragraphs that used to follow "8800010: Canto Coral II (Requisito)",
# while keeping the blank paragraph and the page-break paragraph that
# come right after them.

$d = $word.ActiveDocument

$startMarker = "8800010: Canto Coral II (Requisito)"
$endMarker   = "Original theme under Creative Commons Attribution"

$paras = $d.Paragraphs
$startIndex = -1
$endIndex   = -1

for ($i = 1; $i -le $paras.Count; $i++) {
    $t = $paras.Item($i).Range.Text
    if ($t -like "*$startMarker*") {
        $startIndex = $i
    }
    if ($t -like "*$endMarker*") {
        $endIndex = $i
        break
    }
}

if ($startIndex -gt 0 -and $endIndex -gt $startIndex) {
    # The paragraph right after the "8800010..." one begins the block we
    # want to delete; the block ends with the paragraph containing the
    # copyright/footer text (inclusive), paragraph mark included so the
    # subsequent paragraphs collapse upward cleanly.
    $deleteFrom = $paras.Item($startIndex + 1).Range.Start
    $deleteTo   = $paras.Item($endIndex).Range.End

    $rng = $d.Range($deleteFrom, $deleteTo)
    $rng.Delete()
}
